$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.976.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.781.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5382"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3763"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07438"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.092"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.081"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.216"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.776.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.27%  "

$ws.Range("E18").Value = "  -1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06438"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.877"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.998.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.16%  "

$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.71%  "

$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.981.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.272"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.103"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1053"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.520"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -3.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06417"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02276"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.426"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.444"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6133"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.173"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.664"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5736"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.186"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.920"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06789"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.68%  "
